$d = $word.ActiveDocument

# Locate the start of the block to replace: the paragraph containing "Faire un panier (PHP)"
$startRange = $d.Content
$startFound = $startRange.Find.Execute("Faire un panier (PHP)")
if (-not $startFound) {
    throw "Could not locate start paragraph 'Faire un panier (PHP)'"
}
$startPos = $startRange.Paragraphs(1).Range.Start

# Locate the end of the block to replace: end of the paragraph containing "Mise a jour du MCD (si besoin)"
$endRange = $d.Content
$endFound = $endRange.Find.Execute("Mise a jour du MCD (si besoin)")
if (-not $endFound) {
    throw "Could not locate end paragraph 'Mise a jour du MCD (si besoin)'"
}
$endPos = $endRange.Paragraphs(1).Range.End

# Range covering the whole block of four task-list paragraphs that need replacing
$target = $d.Range($startPos, $endPos)

$openXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Demande d’ajout de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>supplements</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Boissons, viandes, frites, fromage…) (PHP)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Faire un panier (PHP)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Passer la commande (PHP)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Faire des comptes employés</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Faire un compte admin</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Mise </w:t></w:r><w:r><w:t>à</w:t></w:r><w:r><w:t xml:space="preserve"> jour des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mockups</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Mise </w:t></w:r><w:r><w:t>à</w:t></w:r><w:r><w:t xml:space="preserve"> jour du MCD (si besoin)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($openXml)

Write-Host "Replaced task list block."
